$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.134.44"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "3.813.95"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'628.79"
$ws.Range("E5").Value = "  +5.08%  "
$ws.Range("D6").Value = "'165.08"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "3.811.14"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "'35.96"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "4.452.19"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "3.699.30"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "69.138.98"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "'7.12"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'465.76"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "'9.66"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "'0.0000151"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").Value = "'83.61"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "3.961.56"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +3.70%  "
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").Value = "'29.15"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").Value = "'0.149"
$ws.Range("E38").Value = "  +7.68%  "
$ws.Range("D39").Value = "'3.47"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "'0.976"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D44").Value = "'157.00"
$ws.Range("E44").Value = "  +3.54%  "
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("E46").Value = "  +6.08%  "
$ws.Range("D47").Value = "'43.14"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "'46.91"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("D50").Value = "'8.44"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "'0.000278"
$ws.Range("E51").Value = "  +12.38%  "
